# Crowdin update: add English ("英語") translations to column E for the
# Coroner-feature key/value rows (rows 26-46) on the "Main" sheet, matching
# the newly-added shared strings (Coroner, CoronerReport, NoDeadBodies, ...).
# Row 46 already had its English value ("~r~All Units Dismiss~s~") and is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("E26").Value = "Coroner"
$ws.Range("E27").Value = "There is no ~r~dead bodies~s~ nearby you."
$ws.Range("E28").Value = "Requested ~b~{0}~s~ unit to Dispatch."
$ws.Range("E29").Value = "You can check ~b~Coroner's Report~s~ for more information."
$ws.Range("E30").Value = "Have a nice day! Officer!"
$ws.Range("E31").Value = "Press {0} to teleport the backup unit nearby."
$ws.Range("E33").Value = "Coroner Menu"
$ws.Range("E34").Value = "Coroner Report"
$ws.Range("E35").Value = "Report Count: {0}"
$ws.Range("E36").Value = "No Data"
$ws.Range("E38").Value = "Name"
$ws.Range("E39").Value = "Sex"
$ws.Range("E40").Value = "Cause of Death"
$ws.Range("E41").Value = "Died Day"
$ws.Range("E43").Value = "Backup Vehicle"
$ws.Range("E44").Value = "Backup Officer"
